$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the last existing data row (177) down into the
# two new rows so the new date cells pick up the same style (s="1",
# yyyy-mm-dd hh:mm:ss number format) without Excel minting a brand-new
# cellXf entry.
$ws.Range("A177:H177").Copy()
$ws.Range("A178:H179").PasteSpecial(-4122)

$rows = @(
    @{ Row = 178; Date = 45454.2916666667 },
    @{ Row = 179; Date = 45455.2916666667 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 5.15000009536743
    $ws.Cells.Item($row, 4).Value = 5.15000009536743
    $ws.Cells.Item($row, 5).Value = 5.15000009536743
    $ws.Cells.Item($row, 6).Value = 5.15000009536743

    # Column G ("adj_close") stores this value as text (shared string),
    # not as a number, even though it looks numeric. Build the text in an
    # unused scratch cell via a formula (so it is genuinely a string
    # value), then paste just the resulting value into place. This keeps
    # the target cell's text-ness without Excel minting an extra
    # "quote-prefixed"/"@"-formatted style on it.
    $ws.Range("Z1").Formula = "=""5.15000009536743"""
    $ws.Range("Z1").Copy()
    $ws.Cells.Item($row, 7).PasteSpecial(-4163)
    $ws.Range("Z1").Clear()

    $ws.Cells.Item($row, 8).Value = "VLC.MI"
}

$excel.CutCopyMode = 0
